# 4.0.3 model and data
# Split the single "BVTQaZ.csv" and "VTQaZ.csv" rows on the Boolean sheet
# into six per-vehicle-type files each, update view/selection state to
# match the authored workbook, and leave six blank trailing rows.

$wb = $excel.ActiveWorkbook

$wsBoolean = $wb.Worksheets.Item("Boolean")

# --- Row 17 ("trans/BVTQaZ/BVTQaZ.csv") splits into six rows ---------------
$wsBoolean.Rows.Item(18).Resize(5).Insert() | Out-Null

$wsBoolean.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBoolean.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBoolean.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBoolean.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBoolean.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBoolean.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# Rows 18-20 (BVTStL, PVTStL, SRPbVT) have shifted down to 23-25; the old
# VTQaZ.csv row (originally 21) is now at row 26.

# --- Row 26 ("trans/VTQaZ/VTQaZ.csv") splits into six rows ------------------
$wsBoolean.Rows.Item(27).Resize(5).Insert() | Out-Null

$wsBoolean.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBoolean.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBoolean.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBoolean.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBoolean.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBoolean.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# VTStFES.csv (old row 22) is now at row 32.

# --- Six trailing blank (but formatted) rows --------------------------------
$wsBoolean.Rows.Item(33).Resize(6).Insert() | Out-Null

# --- View / selection state --------------------------------------------------
$wsInteger = $wb.Worksheets.Item("Integer")
$wsInteger.Range("A13").Select() | Out-Null

$wsBoolean.Activate()
$wsBoolean.Range("A32").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("A1").Select() | Out-Null
